$wb = $excel.ActiveWorkbook

# Update the "addListItem" sheet: eli_text / eli_code_text value PuneAU -> PuneAV
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "PuneAV"

# Update the "createUser" sheet: increment the use_increment seed 1097 -> 1098
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1098
